$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.180.68"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.268.16"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.47"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.27"
$ws.Range("E6").Value = "  +1.51%  "

$ws.Range("E7").Value = "  -1.17%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -0.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.17"
$ws.Range("E10").Value = "  -2.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  -1.82%  "

$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.87"
$ws.Range("E13").Value = "  +1.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.619.96"
$ws.Range("E14").Value = "  -0.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.70"
$ws.Range("E15").Value = "  +1.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.258.27"
$ws.Range("E16").Value = "  -1.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.791"
$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.098.44"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.28"
$ws.Range("E19").Value = "  -3.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  -1.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.15"

$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("E25").Value = "  +1.42%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -2.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.62"
$ws.Range("E28").Value = "  +4.12%  "

$ws.Range("E29").Value = "  -0.41%  "

$ws.Range("E30").Value = "  +1.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "162.40"
$ws.Range("E31").Value = "  +0.32%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("E34").Value = "  +2.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.62"
$ws.Range("E35").Value = "  +2.04%  "

$ws.Range("E36").Value = "  -2.56%  "

$ws.Range("E37").Value = "  -0.43%  "

$ws.Range("E38").Value = "  -3.82%  "

$ws.Range("E39").Value = "  -1.13%  "

$ws.Range("E40").Value = "  -1.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.09"
$ws.Range("E41").Value = "  -2.31%  "

$ws.Range("E42").Value = "  +3.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.951.32"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.93"
$ws.Range("E44").Value = "  -3.21%  "

$ws.Range("E45").Value = "  -1.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.96"
$ws.Range("E46").Value = "  -1.90%  "

$ws.Range("E47").Value = "  -2.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.83"
$ws.Range("E48").Value = "  +0.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.491.79"
$ws.Range("E49").Value = "  -0.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.08"
$ws.Range("E50").Value = "  -0.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.63"
$ws.Range("E51").Value = "  -1.96%  "
